# Wrote results section in final report.
#
# 1) Fix a few scores on Sheet1 (recalculated totals follow automatically).
# 2) Add a new "Sheet2" after Sheet1 summarising the explore-course results,
#    grouped into three blocks (A9:H21 misc picks, A25:H35 "Zelenski" block,
#    A70:H75 "Sahami" block) each with a SUM row.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# ---- Sheet1 corrections ------------------------------------------------
$ws1.Range("G4").Value  = 100
$ws1.Range("H9").Value  = 15
$ws1.Range("H16").Value = 15

# ---- Add Sheet2 right after Sheet1 -------------------------------------
$ws2 = $wb.Worksheets.Add([Type]::Missing, $ws1)

function Fill-Row($ws, $row, $text, $g, $h) {
    $ws.Range("A$row").Value = $text
    $ws.Range("G$row").Value = $g
    $ws.Range("H$row").Value = $h
}

# Block 1 (rows 9-21, sparse)
Fill-Row $ws2 9  "courses taught by julie zelenski" 0 100
Fill-Row $ws2 10 "courses taught by mehran sahami"  0 200
Fill-Row $ws2 11 "leon simon"                        0 100
Fill-Row $ws2 17 "amelang107b and csre14n"            0 0
Fill-Row $ws2 19 "math52h and cs105"                  0 100
Fill-Row $ws2 21 "amelang129a and cs109"               0 100

# Long-text rows wrap, matching Sheet1's treatment of the same strings.
$ws2.Range("A10").WrapText = $true
$ws2.Rows.Item(10).RowHeight = 48
$ws2.Range("A17").WrapText = $true
$ws2.Rows.Item(17).RowHeight = 36

# Block 2 (rows 25-34) + SUM row 35
Fill-Row $ws2 25 "introduction to computing principles" 100 0
Fill-Row $ws2 26 "the mathematics of the rubik's cube"   0   0
Fill-Row $ws2 27 "third-year persian, second quarter"    0   0
Fill-Row $ws2 28 "identity and popular music (femgen 140g, music 140g) , csre140g" 0 100
Fill-Row $ws2 29 "what is hemispheric"                   100 0
Fill-Row $ws2 30 "first-year hausa"                       0   0
Fill-Row $ws2 31 "first-year hebrew, first quarter (jewishst 101a)" 0 15
Fill-Row $ws2 32 "comparative fictions of ethnicity (amstud 51q, complit 51q) , csre51q" 0 100
Fill-Row $ws2 33 "growing up bilingual (chilatst 14n, educ 114n)" 0 0
Fill-Row $ws2 34 "digital dilemmas"                       100 100
$ws2.Range("G35").Formula = "=SUM(G25:G34)"
$ws2.Range("H35").Formula = "=SUM(H25:H34)"

# Block 3 (rows 70-74) + SUM row 75
Fill-Row $ws2 70 "amelang110a"   100 100
Fill-Row $ws2 71 "math120"       100 100
Fill-Row $ws2 72 "amelang 144b"  100 60
Fill-Row $ws2 73 "math 53"       0   15
Fill-Row $ws2 74 "amelang128b"   100 100
$ws2.Range("G75").Formula = "=SUM(G70:G74)"
$ws2.Range("H75").Formula = "=SUM(H70:H74)"

# New sheet's page setup uses Excel's plain (inch) defaults, not Sheet1's
# custom (cm-derived) margins.
$ws2.PageSetup.LeftMargin   = 0.75 * 72
$ws2.PageSetup.RightMargin  = 0.75 * 72
$ws2.PageSetup.TopMargin    = 1    * 72
$ws2.PageSetup.BottomMargin = 1    * 72
$ws2.PageSetup.HeaderMargin = 0.5  * 72
$ws2.PageSetup.FooterMargin = 0.5  * 72
$ws2.PageSetup.Orientation  = 1

# ---- View state: Sheet2 scrolled down with F60 selected ---------------
[void]$ws2.Range("F60").Select()

# ---- Sheet1 stays the tab that is active/selected on save -------------
[void]$ws1.Range("G5").Select()
$ws1.Activate()
